$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.745.75"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.48%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.468.42"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.84%  "
$ws.Range("E4").Value = "  +0.52%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "92.93"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.72%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.550"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.52%  "
$ws.Range("E8").Value = "  +0.38%  "
$ws.Range("E9").Value = "  +3.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.70"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0841"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +7.20%  "
$ws.Range("E12").Value = "  +0.34%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.847.34"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.94%  "
$ws.Range("E14").Value = "  -0.16%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.77"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.438.90"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.74%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.780"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.82%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "41.715.63"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.09%  "
$ws.Range("E19").Value = "  +2.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0946"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.45%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.02"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.43%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.42"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.65%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "239.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.72"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.10%  "
$ws.Range("E25").Value = "  +1.00%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.61"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.68%  "
$ws.Range("E28").Value = "  +0.73%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.78"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.55"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.79%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "155.84"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.97%  "
$ws.Range("E32").Value = "  +2.32%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.57"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.18%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0762"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.84%  "
$ws.Range("E35").Value = "  +0.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.49"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.35%  "
$ws.Range("E37").Value = "  -1.63%  "
$ws.Range("E38").Value = "  +1.26%  "
$ws.Range("E39").Value = "  +0.87%  "
$ws.Range("E40").Value = "  -1.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.93%  "
$ws.Range("E42").Value = "  +0.41%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.977.18"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.32%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.01"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.84%  "
$ws.Range("E45").Value = "  -0.55%  "
$ws.Range("E46").Value = "  -1.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.04"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.16%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.703.45"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.11%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "96.95"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "67.10"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.40%  "
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "52.46"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.25%  "
